# Re-style the "{qrcode}" placeholder textbox (slide 1, shape "TextBox 8")
# into a small square QR-code placeholder box:
#   - reposition/resize to a 360000x360000 EMU square
#   - add an accent1-colored outline
#   - center the text, zero out the internal margins, switch autofit mode
#   - shrink the font from 6pt to 5pt

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(5)

if ($sh.Name -ne "TextBox 8") {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        if ($s.Shapes.Item($i).Name -eq "TextBox 8") {
            $sh = $s.Shapes.Item($i)
            break
        }
    }
}

# --- position & size: EMU 2748382,1305047 / 360000x360000 -------------------
# (PowerPoint's COM surface works in points; 1 pt = 12700 EMU)
$sh.Left   = 216.40803
$sh.Top    = 102.75961
$sh.Width  = 28.34646
$sh.Height = 28.34646

# --- outline: solid accent1 line --------------------------------------------
$sh.Line.ForeColor.ObjectThemeColor = 5   # msoThemeColorAccent1

# --- text box body: centered, no internal margins, shrink-text-to-fit -------
$tf = $sh.TextFrame
$tf.AutoSize      = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>
$tf.MarginLeft    = 0
$tf.MarginTop     = 0
$tf.MarginRight   = 0
$tf.MarginBottom  = 0
$tf.VerticalAnchor = 3  # msoAnchorMiddle -> anchor="ctr"

# --- font size: 6pt -> 5pt ---------------------------------------------------
$tf.TextRange.Font.Size = 5
